$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: introduce the brand-new item names in the same order the
# author typed them, so the shared-string table grows in that order
# (existing "Barq's" / "Java Monster" entries become unreferenced once
# overwritten below, and get pruned on save).
$ws.Range("A3").Value = "Barq's Root Beer"
$ws.Range("A7").Value = "Java Monster Mean Bean"
$ws.Range("B7").Value = 4
$ws.Range("A9").Value = "Java Monster Irish Crème"
$ws.Range("B9").Value = 4
$ws.Range("A10").Value = "Java Monster Café Latte"
$ws.Range("B10").Value = 4
$ws.Range("A11").Value = "Java Monster Loca Moca"
$ws.Range("B11").Value = 4
$ws.Range("A12").Value = "Java Monster Triple Shot French Vanilla"
$ws.Range("B12").Value = 5
$ws.Range("A13").Value = "Java Monster Triple Shot Mocha"
$ws.Range("B13").Value = 5

# --- Phase 2: lay every row out in its final (alphabetised) order/value.
$final = @(
    @("Item", "Cost"),
    @("7up", 2),
    @("Barq's Root Beer", 2),
    @("Cheerwine", 3),
    @("Coca-Cola", 2),
    @("Dr. Pepper", 2),
    @("Java Monster Café Latte", 4),
    @("Java Monster Irish Crème", 4),
    @("Java Monster Loca Moca", 4),
    @("Java Monster Mean Bean", 4),
    @("Java Monster Triple Shot French Vanilla", 5),
    @("Java Monster Triple Shot Mocha", 5),
    @("Vernors Ginger Ale", 3)
)

$row = 1
foreach ($entry in $final) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

# Widen column A to fit the new, longer item names (target stored width
# 28.42578125 character-units; the host quantises ColumnWidth writes to
# 1/6-character steps, so feed it the input that lands on the nearest
# reachable grid point, 28.5).
$ws.Columns.Item(1).ColumnWidth = 27.65

# Touch page setup so the worksheet prints in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection on the last data cell, matching the saved view.
$ws.Range("B13").Select()
